$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 249.5
$ws.Range("I6").Value = 249.5
$ws.Range("K6").Value = 748.5
$ws.Range("M6").Value = -636.5
$ws.Range("H8").Value = 73
$ws.Range("I8").Value = 73
$ws.Range("K8").Value = 219
$ws.Range("M8").Value = -80
$ws.Range("H18").Value = 5015.8335
$ws.Range("I18").Value = 4960
$ws.Range("J18").Value = 5295
$ws.Range("K18").Value = 4960
$ws.Range("L18").Value = 5295
$ws.Range("M18").Value = -4676
$ws.Range("N18").Value = -5863
$ws.Range("H19").Value = 2063
$ws.Range("I19").Value = 1875.25
$ws.Range("J19").Value = 2250.75
$ws.Range("K19").Value = 1875.25
$ws.Range("L19").Value = 2250.75
$ws.Range("M19").Value = -1700.25
$ws.Range("N19").Value = -2600.75
$ws.Range("H31").Value = 199.5
$ws.Range("I31").Value = 199
$ws.Range("K31").Value = 597
$ws.Range("M31").Value = -367
$ws.Range("H42").Value = 313
$ws.Range("I42").Value = 283.625
$ws.Range("J42").Value = 360
$ws.Range("K42").Value = 850.875
$ws.Range("L42").Value = 1080
$ws.Range("M42").Value = -620.875
$ws.Range("N42").Value = -1540
$ws.Range("H53").Value = 319.57144
$ws.Range("I53").Value = 87
$ws.Range("J53").Value = 901
$ws.Range("K53").Value = 87
$ws.Range("L53").Value = 901
$ws.Range("M53").Value = 550
$ws.Range("N53").Value = -2175
$ws.Range("H80").Value = 465.875
$ws.Range("I80").Value = 390.4
$ws.Range("K80").Value = 1171.2
$ws.Range("M80").Value = -173.1999999999998
$ws.Range("H83").Value = 465.875
$ws.Range("I83").Value = 390.4
$ws.Range("K83").Value = 3513.6
$ws.Range("M83").Value = 1478.4
$ws.Range("H97").Value = 432.66666
$ws.Range("J97").Value = 432.66666
$ws.Range("L97").Value = 1297.99998
$ws.Range("N97").Value = -2289.99998
$ws.Range("H132").Value = 1854.4546
$ws.Range("I132").Value = 1854.4546
$ws.Range("K132").Value = 5563.3638
$ws.Range("M132").Value = -3033.3638
$ws.Range("H137").Value = 3241.24
$ws.Range("I137").Value = 1153.6364
$ws.Range("K137").Value = 3460.9092
$ws.Range("M137").Value = -910.9092000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 58030
$ws.Range("J44").Value = 58030
$ws.Range("L44").Value = 58030
$ws.Range("N44").Value = -59006
$ws.Range("H45").Value = 2550
$ws.Range("I45").Value = 2550
$ws.Range("K45").Value = 2550
$ws.Range("M45").Value = -2173
$ws.Range("H61").Value = 4971.727
$ws.Range("I61").Value = 5169.2
$ws.Range("K61").Value = 5169.2
$ws.Range("M61").Value = -4957.2
$ws.Range("H132").Value = 2409.85
$ws.Range("I132").Value = 1978.0625
$ws.Range("K132").Value = 5934.1875
$ws.Range("M132").Value = -3404.1875
$ws.Range("H136").Value = 4971.727
$ws.Range("I136").Value = 5169.2
$ws.Range("K136").Value = 15507.6
$ws.Range("M136").Value = -12957.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1623.8
$ws.Range("I86").Value = 1037.3334
$ws.Range("J86").Value = 2503.5
$ws.Range("K86").Value = 1037.3334
$ws.Range("L86").Value = 2503.5
$ws.Range("M86").Value = 85.66660000000002
$ws.Range("N86").Value = -4749.5
$ws.Range("H89").Value = 1623.8
$ws.Range("I89").Value = 1037.3334
$ws.Range("J89").Value = 2503.5
$ws.Range("K89").Value = 5186.666999999999
$ws.Range("L89").Value = 12517.5
$ws.Range("M89").Value = 429.3330000000005
$ws.Range("N89").Value = -23749.5
$ws.Range("H134").Value = 4959.222
$ws.Range("I134").Value = 4658.2856
$ws.Range("K134").Value = 13974.8568
$ws.Range("M134").Value = -11439.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1510.7368
$ws.Range("I31").Value = 1450.25
$ws.Range("K31").Value = 1450.25
$ws.Range("M31").Value = -1155.25
$ws.Range("H34").Value = 1510.7368
$ws.Range("I34").Value = 1450.25
$ws.Range("K34").Value = 1450.25
$ws.Range("M34").Value = -1248.25
$ws.Range("H58").Value = 2418.8
$ws.Range("I58").Value = 2447
$ws.Range("J58").Value = 2400
$ws.Range("K58").Value = 2447
$ws.Range("L58").Value = 2400
$ws.Range("M58").Value = -2244
$ws.Range("N58").Value = -2806
$ws.Range("H60").Value = 36666.92
$ws.Range("J60").Value = 36666.92
$ws.Range("L60").Value = 36666.92
$ws.Range("N60").Value = -37688.92
$ws.Range("H134").Value = 4131.4614
$ws.Range("I134").Value = 4191.8
$ws.Range("K134").Value = 12575.4
$ws.Range("M134").Value = -10040.4
$ws.Range("H136").Value = 2418.8
$ws.Range("I136").Value = 2447
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 7341
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -4791
$ws.Range("N136").Value = -12300
$ws.Range("H141").Value = 119666
$ws.Range("J141").Value = 129999.5
$ws.Range("L141").Value = 129999.5
$ws.Range("N141").Value = -140359.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 55393480
$ws.Range("I4").Value = 4288838.5
$ws.Range("J4").Value = 100110040
$ws.Range("K4").Value = 12866515.5
$ws.Range("L4").Value = 300330120
$ws.Range("M4").Value = -12866403.5
$ws.Range("N4").Value = -300330344
$ws.Range("H10").Value = 150
$ws.Range("I10").Value = 150
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 450
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -311
$ws.Range("N10").Value = ""
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3406.1333
$ws.Range("J46").Value = 5916.1665
$ws.Range("L46").Value = 5916.1665
$ws.Range("N46").Value = -6292.1665
$ws.Range("H136").Value = 3108.1667
$ws.Range("I136").Value = 3108.1667
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9324.500100000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6774.500100000001
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 4000
$ws.Range("I7").Value = 4000
$ws.Range("K7").Value = 4000
$ws.Range("M7").Value = -3887
$ws.Range("H132").Value = 2800.52
$ws.Range("I132").Value = 1684
$ws.Range("J132").Value = 4010.0833
$ws.Range("K132").Value = 5052
$ws.Range("L132").Value = 12030.2499
$ws.Range("M132").Value = -2522
$ws.Range("N132").Value = -17090.2499
